# "Added My Timesheets and New Folders"
# This weekly timesheet is rolled forward to the new week (March 8-14,
# 2021), the now-unused "Week Total" label in I5 is cleared, the date
# columns are widened slightly to fit, and the selection is moved back
# to the top of the entry grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Roll the week's dates (row 5, columns B:H) forward to Mar 8-14, 2021.
$ws.Range("B5").Value = 44263
$ws.Range("C5").Value = 44264
$ws.Range("D5").Value = 44265
$ws.Range("E5").Value = 44266
$ws.Range("F5").Value = 44267
$ws.Range("G5").Value = 44268
$ws.Range("H5").Value = 44269

# The "Week Total" label that used to sit in I5 is no longer used.
$ws.Range("I5").Value = ""

# Slightly widen the date columns (D:H) to fit the new values.
$ws.Range("D1:H1").ColumnWidth = 6.42

# Move the active selection back up to H6.
$ws.Range("H6").Select()
